$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (08-09-2021 .. 15-09-2021) appended below the existing data.
# Column A holds the date as TEXT (matching the existing column), so force the
# cell format to Text before assignment, then drop back to the default "Normal"
# style once the literal is stored (keeps the style index untouched like the rest
# of column A, which carries no explicit "s" attribute).
$dateRange = $ws.Range("A174:A179")
$dateRange.NumberFormat = "@"

$ws.Range("A174").Value = "08-09-2021"
$ws.Range("B174").Value = 0.12
$ws.Range("C174").Value = 0.97
$ws.Range("D174").Value = 0.41
$ws.Range("E174").Value = 0.79
$ws.Range("F174").Value = 2.67
$ws.Range("G174").Value = 2.03
$ws.Range("H174").Value = 1.46
$ws.Range("I174").Value = 1.67

$ws.Range("A175").Value = "09-09-2021"
$ws.Range("B175").Value = 0.13
$ws.Range("C175").Value = 0.5
$ws.Range("D175").Value = 0.36
$ws.Range("E175").Value = 1.19
$ws.Range("F175").Value = 3.62
$ws.Range("G175").Value = 1.69
$ws.Range("H175").Value = 0.75
$ws.Range("I175").Value = 3.01

$ws.Range("A176").Value = "10-09-2021"
$ws.Range("B176").Value = 0.13
$ws.Range("C176").Value = 0.5600000000000001
$ws.Range("D176").Value = 0.35
$ws.Range("E176").Value = 0.84
$ws.Range("F176").Value = 1
$ws.Range("G176").Value = 2.4
$ws.Range("H176").Value = 1.87
$ws.Range("I176").Value = 2.41

$ws.Range("A177").Value = "13-09-2021"
$ws.Range("B177").Value = 0.13
$ws.Range("C177").Value = 1.19
$ws.Range("D177").Value = 0.38
$ws.Range("E177").Value = 1.18
$ws.Range("F177").Value = 2.29
$ws.Range("G177").Value = 2.61
$ws.Range("H177").Value = 1.09
$ws.Range("I177").Value = 2.18

$ws.Range("A178").Value = "14-09-2021"
$ws.Range("B178").Value = 0.13
$ws.Range("C178").Value = 0.73
$ws.Range("D178").Value = 0.38
$ws.Range("E178").Value = 0.64
$ws.Range("F178").Value = 2.24
$ws.Range("G178").Value = 3.27
$ws.Range("H178").Value = 0.92
$ws.Range("I178").Value = 2

$ws.Range("A179").Value = "15-09-2021"
$ws.Range("B179").Value = 0.13
$ws.Range("C179").Value = 0.7
$ws.Range("D179").Value = 0.42
$ws.Range("E179").Value = 0.61
$ws.Range("F179").Value = 3.48
$ws.Range("G179").Value = 2.07
$ws.Range("H179").Value = 1.08
$ws.Range("I179").Value = 2.42

$dateRange.Style = "Normal"

